# The workbook's "Recorded By" column (G) contains values that combine the
# recording system name with a user's email address. Rows that currently
# read "System, dnasr281@gmail.com" need to have the two parts swapped to
# read "dnasr281@gmail.com, System". Rows whose "Recorded By" value is only
# "System" or only the email address (no comma-joined pair) must stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
        $changed = $changed + 1
    }
}

Write-Output "Updated $changed 'Recorded By' cells in column G"
